$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.579.93'
$ws.Range('E2').Value = '  +0.33%  '

$ws.Range('D3').Value = '1.912.73'
$ws.Range('E3').Value = '  +0.46%  '

$ws.Range('E4').Value = '  -0.08%  '

$ws.Range('D5').Value = '0.705'
$ws.Range('E5').Value = '  +8.98%  '

$ws.Range('D6').Value = '246.93'
$ws.Range('E6').Value = '  +0.18%  '

$ws.Range('E7').Value = '  -0.04%  '

$ws.Range('D8').Value = "'41.10"
$ws.Range('E8').Value = '  -1.87%  '

$ws.Range('E9').Value = '  +4.24%  '

$ws.Range('D10').Value = '52.65'
$ws.Range('E10').Value = '  +8.05%  '

$ws.Range('E11').Value = '  +3.89%  '

$ws.Range('D12').Value = '0.0989'
$ws.Range('E12').Value = '  -0.93%  '

$ws.Range('D13').Value = '2.189.27'
$ws.Range('E13').Value = '  +0.42%  '

$ws.Range('D14').Value = '12.65'
$ws.Range('E14').Value = '  +2.03%  '

$ws.Range('D15').Value = '0.718'
$ws.Range('E15').Value = '  +3.07%  '

$ws.Range('D16').Value = '4.93'
$ws.Range('E16').Value = '  +2.34%  '

$ws.Range('D17').Value = '1.908.45'
$ws.Range('E17').Value = '  +0.04%  '

$ws.Range('D18').Value = '35.557.45'
$ws.Range('E18').Value = '  +0.17%  '

$ws.Range('D19').Value = '73.41'
$ws.Range('E19').Value = '  +1.98%  '

$ws.Range('E20').Value = '  -0.03%  '

$ws.Range('D21').Value = '13.23'
$ws.Range('E21').Value = '  +4.89%  '

$ws.Range('D22').Value = '243.15'
$ws.Range('E22').Value = '  +0.06%  '

$ws.Range('E23').Value = '  +5.11%  '

$ws.Range('D25').Value = '2.32'
$ws.Range('E25').Value = '  +0.92%  '

$ws.Range('E26').Value = '  +4.20%  '

$ws.Range('D27').Value = '169.56'
$ws.Range('E27').Value = '  -1.27%  '

$ws.Range('E28').Value = '  +1.97%  '

$ws.Range('D29').Value = '18.88'
$ws.Range('E29').Value = '  +5.12%  '

$ws.Range('D30').Value = '0.133'
$ws.Range('E30').Value = '  +4.52%  '

$ws.Range('D31').Value = '4.202.02'
$ws.Range('E31').Value = '  +21.58%  '

$ws.Range('E32').Value = '  +3.14%  '

$ws.Range('D33').Value = '0.0576'
$ws.Range('E33').Value = '  +1.30%  '

$ws.Range('E34').Value = '  +1.28%  '

$ws.Range('B35').Value = 'WEMIXToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D35').Value = '1.87'
$ws.Range('E35').Value = '  +5.89%  '

$ws.Range('B36').Value = 'BinanceUSD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D36').Value = '1.01'
$ws.Range('E36').Value = '  -0.12%  '

$ws.Range('D37').Value = '0.919'
$ws.Range('E37').Value = '  -5.54%  '

$ws.Range('E38').Value = '  +11.07%  '

$ws.Range('D39').Value = '2.03'
$ws.Range('E39').Value = '  +0.30%  '

$ws.Range('D40').Value = '17.28'
$ws.Range('E40').Value = '  +10.47%  '

$ws.Range('D41').Value = '98.12'
$ws.Range('E41').Value = '  +7.14%  '

$ws.Range('E42').Value = '  +1.26%  '

$ws.Range('E43').Value = '  +2.82%  '

$ws.Range('E44').Value = '  +1.54%  '

$ws.Range('D45').Value = '1.359.56'
$ws.Range('E45').Value = '  +0.92%  '

$ws.Range('D46').Value = '2.44'
$ws.Range('E46').Value = '  +2.35%  '

$ws.Range('B47').Value = 'MultiversX'
$ws.Range('C47').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D47').Value = '46.32'
$ws.Range('E47').Value = '  -7.25%  '

$ws.Range('B48').Value = 'HuobiToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D48').Value = '2.42'
$ws.Range('E48').Value = '  +0.36%  '

$ws.Range('B49').Value = 'MXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D49').Value = '2.79'
$ws.Range('E49').Value = '  +1.40%  '

$ws.Range('D50').Value = '12.22'
$ws.Range('E50').Value = '  -4.94%  '

$ws.Range('D51').Value = '6.55'
$ws.Range('E51').Value = '  -0.57%  '
